# Workbook update for EPEX Spot prices data
#  - "Prix Spot" sheet: insert a new date column (06-dec) right after the
#    existing 05-dec column (before the old "01-oct." column, i.e. at EI),
#    shifting every later column one to the right. The new column has no
#    data yet, so every hourly row gets a "-" placeholder like the other
#    not-yet-populated December day columns.
#  - "Gaz" and "CO2" sheets: append the next day's quote (2025-12-04) as a
#    new row at the bottom.

$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert new "06-dec" column at EI -------------------
$ws = $wb.Worksheets.Item("Prix Spot")

$ws.Columns("EI:EI").Insert()

$ws.Range("EI1").Value = "06-dec"
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 139).Value = "-"
}

# --- Sheet "Gaz": append 2025-12-04 row ------------------------------------
# Force the date column to stay plain text (matches every other row in the
# column, which are inline strings rather than real Excel dates) by setting
# a text number format before assigning, then restoring the default "Normal"
# style so the new row doesn't pick up a stray style index.
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A169").NumberFormat = "@"
$wsGaz.Range("A169").Value = "2025-12-04"
$wsGaz.Range("A169").Style = "Normal"
$wsGaz.Range("B169").Value = 25.95

# --- Sheet "CO2": append 2025-12-04 row ------------------------------------
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A169").NumberFormat = "@"
$wsCo2.Range("A169").Value = "2025-12-04"
$wsCo2.Range("A169").Style = "Normal"
$wsCo2.Range("B169").Value = 82.5
